$wb = $excel.ActiveWorkbook

# Sheets to update: "展览" (sheet1) and "全部类型" (sheet4) share identical data.
$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F ("想去人数")
$updates = @{
    3  = 7737
    9  = 5983
    12 = 31
    13 = 1820
    14 = 1337
    17 = 146
    18 = 5539
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
